$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly logged activity entry on row 9
$ws.Range("B9").Value = 6977
$ws.Range("C9").Value = 43923
$ws.Range("D9").Value = 0.83888888888888891
$ws.Range("E9").Value = 0.84166666666666667
$ws.Range("G9").Value = "Updated files from 1.1 to 1.4"

# Move the active selection to C13, matching where the user left off
$ws.Range("C13").Select()
